$wb = $excel.ActiveWorkbook

# --- Sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 459.41666
$ws.Range("I33").Value = 412.55554
$ws.Range("K33").Value = 412.55554
$ws.Range("M33").Value = -183.55554
$ws.Range("H40").Value = 100002340
$ws.Range("J40").Value = 166668400
$ws.Range("L40").Value = 166668400
$ws.Range("N40").Value = -166668750
$ws.Range("H86").Value = 10921.286
$ws.Range("I86").Value = 4999.5
$ws.Range("K86").Value = 4999.5
$ws.Range("M86").Value = -3876.5
$ws.Range("H89").Value = 10921.286
$ws.Range("I89").Value = 4999.5
$ws.Range("K89").Value = 24997.5
$ws.Range("M89").Value = -19381.5
$ws.Range("H96").Value = 968292.4
$ws.Range("I96").Value = 1281.4445
$ws.Range("K96").Value = 3844.3335
$ws.Range("M96").Value = -2471.3335
$ws.Range("H99").Value = 3161.1
$ws.Range("I99").Value = 842.6
$ws.Range("J99").Value = 5479.6
$ws.Range("K99").Value = 2527.8
$ws.Range("L99").Value = 16438.8
$ws.Range("M99").Value = -1029.8
$ws.Range("N99").Value = -19434.8
$ws.Range("H100").Value = 3631.9614
$ws.Range("I100").Value = 2522.8235
$ws.Range("K100").Value = 2522.8235
$ws.Range("M100").Value = -1981.8235
$ws.Range("H113").Value = 6861.4
$ws.Range("I113").Value = 6326.75
$ws.Range("K113").Value = 6326.75
$ws.Range("M113").Value = -3072.75
$ws.Range("H116").Value = 11808.059
$ws.Range("I116").Value = 4466.4546
$ws.Range("K116").Value = 4466.4546
$ws.Range("M116").Value = -1024.4546
$ws.Range("H132").Value = 1774.58
$ws.Range("I132").Value = 1694.6888
$ws.Range("K132").Value = 5084.0664
$ws.Range("M132").Value = -2554.0664
$ws.Range("H137").Value = 2536.5278
$ws.Range("I137").Value = 2303.6453
$ws.Range("K137").Value = 6910.9359
$ws.Range("M137").Value = -4360.9359
$ws.Range("H138").Value = 3302.4036
$ws.Range("I138").Value = 1820.8125
$ws.Range("J138").Value = 5198.84
$ws.Range("K138").Value = 5462.4375
$ws.Range("L138").Value = 15596.52
$ws.Range("M138").Value = -322.4375
$ws.Range("N138").Value = -25876.52
$ws.Range("H141").Value = 3586.4092
$ws.Range("I141").Value = 2918.756
$ws.Range("J141").Value = 12711
$ws.Range("K141").Value = 8756.268
$ws.Range("L141").Value = 38133
$ws.Range("M141").Value = -3576.268
$ws.Range("N141").Value = -48493

# --- Sheet ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 7365.7744
$ws.Range("I32").Value = 5539.241
$ws.Range("J32").Value = 15514.923
$ws.Range("K32").Value = 5539.241
$ws.Range("L32").Value = 15514.923
$ws.Range("M32").Value = -5252.241
$ws.Range("N32").Value = -16088.923
$ws.Range("H141").Value = 99999
$ws.Range("J141").Value = 99999
$ws.Range("L141").Value = 99999
$ws.Range("N141").Value = -110359

# --- Sheet BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H15").Value = 0
$ws.Range("J15").Value = 0
$ws.Range("L15").Value = 0
$ws.Range("N15").Value = ""
$ws.Range("H20").Value = 4802.778
$ws.Range("I20").Value = 6864.0557
$ws.Range("J20").Value = 2741.5
$ws.Range("K20").Value = 6864.0557
$ws.Range("L20").Value = 2741.5
$ws.Range("M20").Value = -6617.0557
$ws.Range("N20").Value = -3235.5
$ws.Range("H86").Value = 2954.639
$ws.Range("I86").Value = 1949.0385
$ws.Range("J86").Value = 5569.2
$ws.Range("K86").Value = 1949.0385
$ws.Range("L86").Value = 5569.2
$ws.Range("M86").Value = -826.0385000000001
$ws.Range("N86").Value = -7815.2
$ws.Range("H89").Value = 2954.639
$ws.Range("I89").Value = 1949.0385
$ws.Range("J89").Value = 5569.2
$ws.Range("K89").Value = 9745.192500000001
$ws.Range("L89").Value = 27846
$ws.Range("M89").Value = -4129.192500000001
$ws.Range("N89").Value = -39078
$ws.Range("H105").Value = 957223.9
$ws.Range("I105").Value = 1635948.1
$ws.Range("J105").Value = 7009.9
$ws.Range("K105").Value = 1635948.1
$ws.Range("L105").Value = 7009.9
$ws.Range("M105").Value = -1634201.1
$ws.Range("N105").Value = -10503.9
$ws.Range("H107").Value = 4218.4546
$ws.Range("I107").Value = 2788.0667
$ws.Range("K107").Value = 2788.0667
$ws.Range("M107").Value = -868.0666999999999

# --- Sheet CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H12").Value = 2667.6667
$ws.Range("I12").Value = 2667.6667
$ws.Range("J12").Value = 0
$ws.Range("K12").Value = 2667.6667
$ws.Range("L12").Value = 0
$ws.Range("M12").Value = -2497.6667
$ws.Range("N12").Value = ""
$ws.Range("H31").Value = 25004360
$ws.Range("I31").Value = 55558932
$ws.Range("K31").Value = 55558932
$ws.Range("M31").Value = -55558637
$ws.Range("H34").Value = 25004360
$ws.Range("I34").Value = 55558932
$ws.Range("K34").Value = 55558932
$ws.Range("M34").Value = -55558730
$ws.Range("H60").Value = 32794.75
$ws.Range("I60").Value = 37059.668
$ws.Range("J60").Value = 20000
$ws.Range("K60").Value = 37059.668
$ws.Range("L60").Value = 20000
$ws.Range("M60").Value = -36548.668
$ws.Range("N60").Value = -21022
$ws.Range("H70").Value = 0
$ws.Range("J70").Value = 0
$ws.Range("L70").Value = 0
$ws.Range("N70").Value = ""
$ws.Range("H73").Value = 0
$ws.Range("J73").Value = 0
$ws.Range("L73").Value = 0
$ws.Range("N73").Value = ""
$ws.Range("H132").Value = 2364.3333
$ws.Range("I132").Value = 2329.913
$ws.Range("K132").Value = 6989.739
$ws.Range("M132").Value = -4459.739
$ws.Range("H134").Value = 2395.5293
$ws.Range("I134").Value = 2085.0833
$ws.Range("J134").Value = 3140.6
$ws.Range("K134").Value = 6255.249899999999
$ws.Range("L134").Value = 9421.799999999999
$ws.Range("M134").Value = -3720.249899999999
$ws.Range("N134").Value = -14491.8

# --- Sheet CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 113.82353
$ws.Range("J2").Value = 104.666664
$ws.Range("L2").Value = 627.999984
$ws.Range("N2").Value = -853.999984
$ws.Range("H80").Value = 27782946
$ws.Range("I80").Value = 55557556
$ws.Range("J80").Value = 8333.333000000001
$ws.Range("K80").Value = 166672668
$ws.Range("L80").Value = 24999.999
$ws.Range("M80").Value = -166671732
$ws.Range("N80").Value = -26871.999
$ws.Range("H83").Value = 27782946
$ws.Range("I83").Value = 55557556
$ws.Range("J83").Value = 8333.333000000001
$ws.Range("K83").Value = 500018004
$ws.Range("L83").Value = 74999.997
$ws.Range("M83").Value = -500013324
$ws.Range("N83").Value = -84359.997
$ws.Range("H133").Value = 31516.934
$ws.Range("I133").Value = 28285.334
$ws.Range("K133").Value = 84856.00199999999
$ws.Range("M133").Value = -79796.00199999999
$ws.Range("H140").Value = 2983.5
$ws.Range("J140").Value = 11592.333
$ws.Range("L140").Value = 34776.999
$ws.Range("N140").Value = -45136.999

# --- Sheet GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 1415.3871
$ws.Range("I102").Value = 1340.5862
$ws.Range("K102").Value = 1340.5862
$ws.Range("M102").Value = 281.4138
$ws.Range("H132").Value = 2166845.2
$ws.Range("I132").Value = 2356.6924
$ws.Range("K132").Value = 7070.0772
$ws.Range("M132").Value = -4540.0772
$ws.Range("H136").Value = 11665.667
$ws.Range("J136").Value = 11665.667
$ws.Range("L136").Value = 34997.001
$ws.Range("N136").Value = -40097.001

# --- Sheet LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 3958.4333
$ws.Range("I132").Value = 2379.1875
$ws.Range("J132").Value = 5763.2856
$ws.Range("K132").Value = 7137.5625
$ws.Range("L132").Value = 17289.8568
$ws.Range("M132").Value = -4607.5625
$ws.Range("N132").Value = -22349.8568
$ws.Range("H136").Value = 6445.6665
$ws.Range("I136").Value = 2448.8572
$ws.Range("K136").Value = 7346.571599999999
$ws.Range("M136").Value = -4796.571599999999
$ws.Range("H140").Value = 277759.34
$ws.Range("J140").Value = 277759.34
$ws.Range("L140").Value = 277759.34
$ws.Range("N140").Value = -288119.34

# --- Sheet WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H48").Value = 0
$ws.Range("I48").Value = 0
$ws.Range("J48").Value = 0
$ws.Range("K48").Value = 0
$ws.Range("L48").Value = 0
$ws.Range("M48").Value = ""
$ws.Range("N48").Value = ""
$ws.Range("H107").Value = 3896.5789
$ws.Range("I107").Value = 2275.0476
$ws.Range("K107").Value = 6825.1428
$ws.Range("M107").Value = -4905.1428
$ws.Range("H122").Value = 2056.1914
$ws.Range("I122").Value = 1872.0883
$ws.Range("K122").Value = 5616.2649
$ws.Range("M122").Value = -3166.2649
$ws.Range("H128").Value = 71785.164
$ws.Range("J128").Value = 71785.164
$ws.Range("L128").Value = 71785.164
$ws.Range("N128").Value = -81745.164
$ws.Range("H136").Value = 3098.4614
$ws.Range("I136").Value = 1173.5
$ws.Range("K136").Value = 3520.5
$ws.Range("M136").Value = -970.5
